$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for columns D, J, K, L, M, P on rows 2-4
$cols = @("D", "J", "K", "L", "M", "P")
$rows = @(2, 3, 4)

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Rotate the values: new row2 <- old row4, new row3 <- old row2, new row4 <- old row3
$mapping = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $orig[$src][$c]
    }
}
